$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A1").Value = 1
$ws.Range("B1").Value = 2
$ws.Range("B1").Select()
